# Regenerate save_data to use K (strikeouts) values instead of Strike# values.
# Only column G ("K") on rows 2-25 changes; all other data is left intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 2
    4  = 3
    5  = 1
    6  = 2
    7  = 3
    8  = 3
    9  = 2
    10 = 0
    11 = 1
    12 = 1
    13 = 1
    14 = 1
    15 = 2
    16 = 1
    17 = 1
    18 = 1
    19 = 1
    20 = 0
    21 = 2
    22 = 1
    23 = 0
    24 = 1
    25 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
